$wb = $excel.ActiveWorkbook

# Source sheets to duplicate: "Raw data_discovery - (4)" and "Raw data_outliers - (4)"
$srcDiscovery = $wb.Worksheets.Item("Raw data_discovery - (4)")
$srcOutliers  = $wb.Worksheets.Item("Raw data_outliers - (4)")

# Copy the discovery sheet to the end of the workbook and rename it
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcDiscovery.Copy($null, $lastSheet)
$newDiscovery = $wb.Worksheets.Item($wb.Worksheets.Count)
$newDiscovery.Name = "Raw data_discovery - (5)"

# Copy the outliers sheet right after the new discovery sheet and rename it
$srcOutliers.Copy($null, $newDiscovery)
$newOutliers = $wb.Worksheets.Item($wb.Worksheets.Count)
$newOutliers.Name = "Raw data_outliers - (5)"
